$wb = $excel.ActiveWorkbook

$wsLinear = $wb.Worksheets.Item("linear")
$wsLinear.Range("B2").Value = 0.6869564407548443
$wsLinear.Range("B3").Value = -1.216704673061194
$wsLinear.Range("B4").Value = 342.0670551671182

$wsNonLinear = $wb.Worksheets.Item("non-linear")
$wsNonLinear.Range("B2").Value = -0.1611150080683665
$wsNonLinear.Range("B3").Value = -3.240443830622895
$wsNonLinear.Range("B4").Value = 499.3887462084699
$wsNonLinear.Range("B5").Value = 0.4311830877526545
$wsNonLinear.Range("B6").Value = -0.5191012567851773
$wsNonLinear.Range("B7").Value = 226.9848487197836
